$d = $word.ActiveDocument

# 1. Update the "Created" date from 2020-05-28 to 2020-06-03
$d.Content.Find.Execute("2020-05-28", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2020-06-03", 2)

# 2. The second R code chunk ("## [1] "Constructing panel for 14 indicators"")
#    is hidden, mirroring the first code chunk: drop the leftover
#    "SourceCode" paragraph style, mark the paragraph mark and the run as
#    vanish (hidden), and collapse the run text down to a single space.
$findRange = $d.Content
$findRange.Find.Execute("Constructing panel for 14 indicators", $true, $false, $false,
                         $false, $false, $true, 1, $false, $null, 0)
$chunkPara = $findRange.Paragraphs(1)
$chunkRange = $chunkPara.Range
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="redoc-codechunk-2"/><w:rPr><w:vanish/></w:rPr></w:pPr><w:r><w:rPr><w:vanish/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$chunkRange.InsertXML($xml)

# 3. Fix the REMITTANCES heading casing
$d.Content.Find.Execute("REMITTANCES", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Remittances", 2)
